# Fruta / hortaliza, semanal
# Insert this week's new price record for Albahaca (Vega Modelo de Temuco)
# at the top of the data block (row 427), pushing the existing history
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(427).Insert()

$ws.Range("A427").Value = 10
$ws.Range("B427").Value = "Vega Modelo de Temuco"
$ws.Range("C427").Value = "La Araucanía"
$ws.Range("D427").Value = 45265
$ws.Range("E427").Value = 9
$ws.Range("F427").Value = 100112052
$ws.Range("G427").Value = "Albahaca"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 50
$ws.Range("K427").Value = 9000
$ws.Range("L427").Value = 9000
$ws.Range("M427").Value = 9000
$ws.Range("N427").Value = "$/paquete"
$ws.Range("O427").Value = "Región Metropolitana"
$ws.Range("P427").Value = 9000
$ws.Range("Q427").Value = 1
$ws.Range("R427").Value = "Hortaliza"
